$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 17.05385513254113; C = 3.792633579969277; D = 13.52141097443091; E = 13.51990321559816; G = 66.04962284693244; H = 23.56213596488568; I = 37.64272558579754; J = 8.194841690971835; K = 14.61160874003793; L = 13.28180907251838 }
    3 = @{ B = 16.96797328073792; C = 3.600829032618984; D = 13.51774323250369; E = 13.53855695207581; G = 65.9097632420541; H = 23.57497338464936; I = 37.64739958468684; J = 8.202168788142011; K = 14.55800253830981; L = 13.29387726346508 }
    4 = @{ B = 16.91947692858822; C = 3.476845375807701; D = 13.51784417202443; E = 13.55164130923285; G = 65.83519696042465; H = 23.58588113656917; I = 37.65530467762488; J = 8.206920049372266; K = 14.52859598485773; L = 13.30326923474735 }
    5 = @{ B = 16.90079771668995; C = 3.424781937901391; D = 13.51847812263653; E = 13.55738374226754; G = 65.8076666362884; H = 23.59108682391552; I = 37.65979219555511; J = 8.208919893642214; K = 14.51750502685206; L = 13.30759530863419 }
    6 = @{ B = 16.89776199170627; C = 3.416044675794668; D = 13.5186192226615; E = 13.55836206877732; G = 65.80326806480065; H = 23.59199716189938; I = 37.66061380740355; J = 8.209255817466616; K = 14.51571756520667; L = 13.30834378307564 }
    7 = @{ B = 16.91922060426366; C = 3.476149426894621; D = 13.51785032008693; E = 13.55171709134864; G = 65.83481409817399; H = 23.5859482625966; I = 37.65536007174239; J = 8.206946761920811; K = 14.52844278194616; L = 13.30332555778809 }
    8 = @{ B = 17.02337552501358; C = 3.727802023352567; D = 13.51965874771927; E = 13.52599675264553; G = 65.9990604561123; H = 23.56593433983102; I = 37.64329220206152; J = 8.197315810326897; K = 14.59240380915839; L = 13.28555905614513 }
    9 = @{ B = 17.26034717691647; C = 4.17116717871981; D = 13.54181455247464; E = 13.48848372174799; G = 66.4102283917596; H = 23.55069428458429; I = 37.6595707077254; J = 8.180423085445168; K = 14.74515980025824; L = 13.26642698461285 }
    10 = @{ B = 17.45311909935694; C = 4.465528290717073; D = 13.56934246114938; E = 13.4687799354629; G = 66.7656247502827; H = 23.55412888701071; I = 37.69585471992483; J = 8.169214657444671; K = 14.87330543500975; L = 13.26191751644773 }
    11 = @{ B = 17.54457010862334; C = 4.592514417974717; D = 13.58428257668566; E = 13.4615168813653; G = 66.9386203483209; H = 23.55886424356766; I = 37.71762949927664; J = 8.164374110194768; K = 14.93487959507419; L = 13.26193070828452 }
    12 = @{ B = 17.57971377724098; C = 4.639599096453145; D = 13.59028494011341; E = 13.45901052193648; G = 67.00573175544567; H = 23.56111290522131; I = 37.72663069039211; J = 8.162578047855598; K = 14.95865141444407; L = 13.26223165347865 }
    13 = @{ B = 17.57212258012322; C = 4.629503234912359; D = 13.58897693293734; E = 13.45953946765948; G = 66.99120732338054; H = 23.56060837235314; I = 37.72465855769982; J = 8.162963221669859; K = 14.953511763854; L = 13.26215369323942 }
    14 = @{ B = 17.54745127428788; C = 4.596408222482347; D = 13.58476950000632; E = 13.4613057945245; G = 66.94410970357607; H = 23.55904011887799; I = 37.71835490393786; J = 8.164225607625523; K = 14.93682629575158; L = 13.26194954172481 }
    15 = @{ B = 17.53240539082049; C = 4.576005876431241; D = 13.58223715040909; E = 13.46241948138427; G = 66.91546881552303; H = 23.55813880946215; I = 37.7145920591681; J = 8.165003661730323; K = 14.9266647146606; L = 13.26186300467588 }
    16 = @{ B = 17.44721587486875; C = 4.457089545115733; D = 13.56841450596313; E = 13.46928876577071; G = 66.7545447213418; H = 23.5538832057168; I = 37.69453754121698; J = 8.169536178909748; K = 14.8693460169562; L = 13.26195812175929 }
    17 = @{ B = 17.39589831484617; C = 4.3823601139119; D = 13.56055206804652; E = 13.47393799462478; G = 66.65870509447822; H = 23.55208470300903; I = 37.68358288990468; J = 8.172382736882277; K = 14.83501136341948; L = 13.26254463342803 }
    18 = @{ B = 17.36673723642336; C = 4.338726644903647; D = 13.55625753546683; E = 13.47677218066922; G = 66.60464828270771; H = 23.55134896253516; I = 37.67777816158916; J = 8.174044317823732; K = 14.81557276314166; L = 13.26307637006457 }
    19 = @{ B = 17.35692565588895; C = 4.323841527434993; D = 13.55484267377527; E = 13.47775929448874; G = 66.58652971351195; H = 23.55115117795664; I = 37.67589804313264; J = 8.174611082530555; K = 14.80904488719238; L = 13.26328982020585 }
    20 = @{ B = 17.40132456604786; C = 4.390382630964633; D = 13.56136548757218; E = 13.47342651284649; G = 66.66879707902977; H = 23.55224524623846; I = 37.684697700943; J = 8.172077200642608; K = 14.83863440248938; L = 13.26246208593503 }
    21 = @{ B = 17.55468412960818; C = 4.606156276686857; D = 13.58599598899886; E = 13.46078036356852; G = 66.95790016894598; H = 23.55948839838236; I = 37.72018595336069; J = 8.163853812977186; K = 14.94171500954109; L = 13.26200148201445 }
    22 = @{ B = 17.65789129070478; C = 4.741333162162849; D = 13.60410210083522; E = 13.45393733942394; G = 67.15616796454241; H = 23.56687663272654; I = 37.74778223915239; J = 8.158694629614063; K = 15.01172810936613; L = 13.26342495444768 }
    23 = @{ B = 17.60254444502209; C = 4.669723274157053; D = 13.5942557256945; E = 13.45745965634144; G = 67.04950509266537; H = 23.56269082165; I = 37.73265155110509; J = 8.161428545629974; K = 14.9741245723803; L = 13.2625077675557 }
    24 = @{ B = 17.39887029027728; C = 4.386757739049631; D = 13.5609970374153; E = 13.47365725134604; G = 66.66423124280739; H = 23.55217173549223; I = 37.68419215817244; J = 8.172215255512995; K = 14.83699548812928; L = 13.26249879960021 }
    25 = @{ B = 17.19286357738084; C = 4.056688749617332; D = 13.53383720545311; E = 13.49725046209204; G = 66.28955317800106; H = 23.55224668036674; I = 37.65089180703875; J = 8.184780906707525; K = 14.70098247045745; L = 13.26992363374515 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value2 = $data[$row][$col]
    }
}

Write-Output "Applied loading_percent updates for case with 380 kV"
